# Update the division-practice table values per the commit's regenerated
# output. Cells are addressed directly via Table.Cell(row, col).Range.Text
# (1-indexed) so duplicate old/new text values elsewhere in the table can't
# cause cross-talk, and the run/paragraph formatting (rFonts, sz, jc) is left
# untouched since we only assign the cell Range's .Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "70÷2=35, 0"
$t.Cell(1, 2).Range.Text = "76÷4=19, 0"
$t.Cell(1, 3).Range.Text = "79÷6=13, 1"
$t.Cell(1, 4).Range.Text = "49÷5=9, 4"
$t.Cell(1, 5).Range.Text = "41÷8=5, 1"

$t.Cell(5, 1).Range.Text = "66÷6=11, 0"
$t.Cell(5, 2).Range.Text = "82÷8=10, 2"
# Cell(5, 3) "67÷2=33, 1" is unchanged.
$t.Cell(5, 4).Range.Text = "64÷6=10, 4"
$t.Cell(5, 5).Range.Text = "38÷3=12, 2"

$t.Cell(9, 1).Range.Text = "82÷8=10, 2"
$t.Cell(9, 2).Range.Text = "33÷3=11, 0"
$t.Cell(9, 3).Range.Text = "68÷6=11, 2"
$t.Cell(9, 4).Range.Text = "82÷9=9, 1"
$t.Cell(9, 5).Range.Text = "59÷9=6, 5"

$t.Cell(13, 1).Range.Text = "46÷2=23, 0"
$t.Cell(13, 2).Range.Text = "25÷7=3, 4"
$t.Cell(13, 3).Range.Text = "61÷5=12, 1"
$t.Cell(13, 4).Range.Text = "22÷8=2, 6"
$t.Cell(13, 5).Range.Text = "27÷3=9, 0"

$t.Cell(17, 1).Range.Text = "97÷7=13, 6"
$t.Cell(17, 2).Range.Text = "47÷9=5, 2"
$t.Cell(17, 3).Range.Text = "79÷6=13, 1"
$t.Cell(17, 4).Range.Text = "20÷4=5, 0"
$t.Cell(17, 5).Range.Text = "67÷8=8, 3"
